$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.4360334873199463
$ws.Range("E2").Value = 62.02870150422677
$ws.Range("F2").Value = 0.002200436538785911
$ws.Range("G2").Value = 0.001863035246001035
$ws.Range("H2").Value = 0.001847652463935592
$ws.Range("I2").Value = 0.001660933182300399
$ws.Range("J2").Value = 0.001660933182300399
$ws.Range("K2").Value = 0.001566895311480357
$ws.Range("L2").Value = 0.001462627750947911
$ws.Range("M2").Value = 0.001462627750947911
$ws.Range("N2").Value = 0.001462627750947911
$ws.Range("O2").Value = 0.00141115928044901
$ws.Range("P2").Value = 0.00141115928044901
$ws.Range("Q2").Value = 0.00141115928044901
$ws.Range("R2").Value = 0.00141115928044901
$ws.Range("S2").Value = 0.00137864408452327
$ws.Range("T2").Value = 0.001327689516936438
$ws.Range("U2").Value = 0.001304081051696565
$ws.Range("V2").Value = 0.001239204944306595
$ws.Range("W2").Value = 0.001239204944306595
$ws.Range("X2").Value = 0.001209136481563874
$ws.Range("Y2").Value = 0.001209136481563874

$ws.Range("C3").Value = 0.5779569149017334
$ws.Range("E3").Value = 67.41487830349433
$ws.Range("F3").Value = 0.002200436538785911
$ws.Range("G3").Value = 0.001710634994383353
$ws.Range("H3").Value = 0.001672910308660677
$ws.Range("I3").Value = 0.001672910308660677
$ws.Range("J3").Value = 0.001672910308660677
$ws.Range("K3").Value = 0.001672910308660677
$ws.Range("L3").Value = 0.001672910308660677
$ws.Range("M3").Value = 0.001637278453653828
$ws.Range("N3").Value = 0.001566499457487587
$ws.Range("O3").Value = 0.001566499457487587
$ws.Range("P3").Value = 0.001476623050347181
$ws.Range("Q3").Value = 0.001476623050347181
$ws.Range("R3").Value = 0.00145515452410499
$ws.Range("S3").Value = 0.00145515452410499
$ws.Range("T3").Value = 0.00145135957306967
$ws.Range("U3").Value = 0.001357455762556827
$ws.Range("V3").Value = 0.001357455762556827
$ws.Range("W3").Value = 0.001339539145054813
$ws.Range("X3").Value = 0.001328385329446757
$ws.Range("Y3").Value = 0.001314130181354665

$ws.Range("C4").Value = 0.4979994297027588
$ws.Range("E4").Value = 64.0498300499803
$ws.Range("F4").Value = 0.002200436538785911
$ws.Range("G4").Value = 0.001906970634048257
$ws.Range("H4").Value = 0.001906970634048257
$ws.Range("I4").Value = 0.001878189102766445
$ws.Range("J4").Value = 0.001811872764294491
$ws.Range("K4").Value = 0.001790154379744452
$ws.Range("L4").Value = 0.001786794190051133
$ws.Range("M4").Value = 0.001653691015456709
$ws.Range("N4").Value = 0.001653691015456709
$ws.Range("O4").Value = 0.00151359276661117
$ws.Range("P4").Value = 0.001484191793877588
$ws.Range("Q4").Value = 0.001481585413819286
$ws.Range("R4").Value = 0.001464963837664854
$ws.Range("S4").Value = 0.001464963837664854
$ws.Range("T4").Value = 0.001464963837664854
$ws.Range("U4").Value = 0.00140250306585408
$ws.Range("V4").Value = 0.001369715229077725
$ws.Range("W4").Value = 0.001277600949566976
$ws.Range("X4").Value = 0.001277600949566976
$ws.Range("Y4").Value = 0.001248534698830025

$ws.Range("C5").Value = 0.4749960899353027
$ws.Range("E5").Value = 64.8660425643684
$ws.Range("F5").Value = 0.002200436538785911
$ws.Range("G5").Value = 0.001624035047760429
$ws.Range("H5").Value = 0.001624035047760429
$ws.Range("I5").Value = 0.001624035047760429
$ws.Range("J5").Value = 0.001624035047760429
$ws.Range("K5").Value = 0.001624035047760429
$ws.Range("L5").Value = 0.001624035047760429
$ws.Range("M5").Value = 0.001624035047760429
$ws.Range("N5").Value = 0.001534112000041226
$ws.Range("O5").Value = 0.001534112000041226
$ws.Range("P5").Value = 0.001534112000041226
$ws.Range("Q5").Value = 0.001447131481633872
$ws.Range("R5").Value = 0.001447131481633872
$ws.Range("S5").Value = 0.00135974801716568
$ws.Range("T5").Value = 0.001349796083284541
$ws.Range("U5").Value = 0.001349796083284541
$ws.Range("V5").Value = 0.001349796083284541
$ws.Range("W5").Value = 0.001293348454267769
$ws.Range("X5").Value = 0.001293348454267769
$ws.Range("Y5").Value = 0.001264445274159228

$ws.Range("C6").Value = 0.479008674621582
$ws.Range("E6").Value = 62.44841208449543
$ws.Range("F6").Value = 0.002128667048956864
$ws.Range("G6").Value = 0.001956752978130786
$ws.Range("H6").Value = 0.001780883019066953
$ws.Range("I6").Value = 0.001756214168716552
$ws.Range("J6").Value = 0.001570473627546315
$ws.Range("K6").Value = 0.001570473627546315
$ws.Range("L6").Value = 0.001570473627546315
$ws.Range("M6").Value = 0.001570473627546315
$ws.Range("N6").Value = 0.001455636820130351
$ws.Range("O6").Value = 0.001455636820130351
$ws.Range("P6").Value = 0.001455636820130351
$ws.Range("Q6").Value = 0.001455636820130351
$ws.Range("R6").Value = 0.001382906810097117
$ws.Range("S6").Value = 0.001382906810097117
$ws.Range("T6").Value = 0.001250094298036282
$ws.Range("U6").Value = 0.001250094298036282
$ws.Range("V6").Value = 0.001250094298036282
$ws.Range("W6").Value = 0.001250094298036282
$ws.Range("X6").Value = 0.001240471153620556
$ws.Range("Y6").Value = 0.001217317974356636

$ws.Range("C7").Value = 0.4480364322662354
$ws.Range("E7").Value = 63.93471908705942
$ws.Range("F7").Value = 0.002083323006429221
$ws.Range("G7").Value = 0.001861054796401135
$ws.Range("H7").Value = 0.001861054796401135
$ws.Range("I7").Value = 0.001819897513340155
$ws.Range("J7").Value = 0.00173151645674204
$ws.Range("K7").Value = 0.00173151645674204
$ws.Range("L7").Value = 0.001554673078201154
$ws.Range("M7").Value = 0.001554673078201154
$ws.Range("N7").Value = 0.001554673078201154
$ws.Range("O7").Value = 0.001554673078201154
$ws.Range("P7").Value = 0.00147869349791234
$ws.Range("Q7").Value = 0.00147869349791234
$ws.Range("R7").Value = 0.001358312776229196
$ws.Range("S7").Value = 0.001358312776229196
$ws.Range("T7").Value = 0.001358312776229196
$ws.Range("U7").Value = 0.00134575518108175
$ws.Range("V7").Value = 0.001314356959203854
$ws.Range("W7").Value = 0.00129968634253707
$ws.Range("X7").Value = 0.001259159729248232
$ws.Range("Y7").Value = 0.001246290820410515

$ws.Range("C8").Value = 0.8669610023498535
$ws.Range("E8").Value = 66.29549042661893
$ws.Range("F8").Value = 0.002090109990741752
$ws.Range("G8").Value = 0.001887526077161699
$ws.Range("H8").Value = 0.001845462964590224
$ws.Range("I8").Value = 0.001790598981263856
$ws.Range("J8").Value = 0.001790598981263856
$ws.Range("K8").Value = 0.001617424259157029
$ws.Range("L8").Value = 0.001617424259157029
$ws.Range("M8").Value = 0.001560051449261245
$ws.Range("N8").Value = 0.001560051449261245
$ws.Range("O8").Value = 0.001394859472989641
$ws.Range("P8").Value = 0.001394859472989641
$ws.Range("Q8").Value = 0.001394859472989641
$ws.Range("R8").Value = 0.001394859472989641
$ws.Range("S8").Value = 0.001394859472989641
$ws.Range("T8").Value = 0.001381592628275572
$ws.Range("U8").Value = 0.001381592628275572
$ws.Range("V8").Value = 0.001354964783069172
$ws.Range("W8").Value = 0.001348880453541368
$ws.Range("X8").Value = 0.001325744764843242
$ws.Range("Y8").Value = 0.001292309754904852

$ws.Range("C9").Value = 0.5470013618469238
$ws.Range("E9").Value = 64.1362535840417
$ws.Range("F9").Value = 0.002200436538785911
$ws.Range("G9").Value = 0.001837855437458066
$ws.Range("H9").Value = 0.001744083847164834
$ws.Range("I9").Value = 0.001744083847164834
$ws.Range("J9").Value = 0.001605822255938871
$ws.Range("K9").Value = 0.001605822255938871
$ws.Range("L9").Value = 0.001590036920401987
$ws.Range("M9").Value = 0.001590036920401987
$ws.Range("N9").Value = 0.001557637483095208
$ws.Range("O9").Value = 0.001557637483095208
$ws.Range("P9").Value = 0.001557637483095208
$ws.Range("Q9").Value = 0.001438178020783093
$ws.Range("R9").Value = 0.001433525836672659
$ws.Range("S9").Value = 0.001394189414677172
$ws.Range("T9").Value = 0.001394189414677172
$ws.Range("U9").Value = 0.001374897254811134
$ws.Range("V9").Value = 0.001347483394345201
$ws.Range("W9").Value = 0.001305683340708224
$ws.Range("X9").Value = 0.001256866154740494
$ws.Range("Y9").Value = 0.001250219368109974

$ws.Range("C10").Value = 0.6820008754730225
$ws.Range("E10").Value = 67.10632048071784
$ws.Range("F10").Value = 0.001966018487533238
$ws.Range("G10").Value = 0.001937019300936586
$ws.Range("H10").Value = 0.001846282152210586
$ws.Range("I10").Value = 0.001846282152210586
$ws.Range("J10").Value = 0.001846282152210586
$ws.Range("K10").Value = 0.001779646611438981
$ws.Range("L10").Value = 0.001601695395008088
$ws.Range("M10").Value = 0.001601695395008088
$ws.Range("N10").Value = 0.00154636388590751
$ws.Range("O10").Value = 0.00154636388590751
$ws.Range("P10").Value = 0.00154636388590751
$ws.Range("Q10").Value = 0.00154636388590751
$ws.Range("R10").Value = 0.001531325767325123
$ws.Range("S10").Value = 0.001528657121966105
$ws.Range("T10").Value = 0.001501542701289172
$ws.Range("U10").Value = 0.001367730075974987
$ws.Range("V10").Value = 0.001367730075974987
$ws.Range("W10").Value = 0.001308115408980854
$ws.Range("X10").Value = 0.001308115408980854
$ws.Range("Y10").Value = 0.001308115408980854

$ws.Range("C11").Value = 0.9950058460235596
$ws.Range("E11").Value = 65.72871479378773
$ws.Range("F11").Value = 0.002200436538785911
$ws.Range("G11").Value = 0.001930433371228211
$ws.Range("H11").Value = 0.001689530768166271
$ws.Range("I11").Value = 0.001689530768166271
$ws.Range("J11").Value = 0.001476023078952277
$ws.Range("K11").Value = 0.001476023078952277
$ws.Range("L11").Value = 0.001476023078952277
$ws.Range("M11").Value = 0.001476023078952277
$ws.Range("N11").Value = 0.001476023078952277
$ws.Range("O11").Value = 0.001476023078952277
$ws.Range("P11").Value = 0.001476023078952277
$ws.Range("Q11").Value = 0.001466575737196787
$ws.Range("R11").Value = 0.00138487360516039
$ws.Range("S11").Value = 0.00138487360516039
$ws.Range("T11").Value = 0.001329246993103893
$ws.Range("U11").Value = 0.001329246993103893
$ws.Range("V11").Value = 0.001329246993103893
$ws.Range("W11").Value = 0.001281261496954926
$ws.Range("X11").Value = 0.001281261496954926
$ws.Range("Y11").Value = 0.001281261496954926
